$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.380.68'
$ws.Range("E2").Value = '  +2.51%  '
$ws.Range("D3").Value = '3.201.21'
$ws.Range("E3").Value = '  +1.85%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''598.17'
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("D6").Value = '''153.65'
$ws.Range("E6").Value = '  +5.60%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.196.09'
$ws.Range("E8").Value = '  +1.90%  '
$ws.Range("E9").Value = '  +2.43%  '
$ws.Range("E10").Value = '  +4.01%  '
$ws.Range("D11").Value = '''6.11'
$ws.Range("E11").Value = '  +6.05%  '
$ws.Range("D12").Value = '''0.473'
$ws.Range("E12").Value = '  +2.95%  '
$ws.Range("D13").Value = '''0.0000255'
$ws.Range("E13").Value = '  +3.22%  '
$ws.Range("D14").Value = '''39.32'
$ws.Range("E14").Value = '  +6.47%  '
$ws.Range("D15").Value = '3.727.45'
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").Value = '''7.43'
$ws.Range("E17").Value = '  +4.64%  '
$ws.Range("D18").Value = '65.083.56'
$ws.Range("E18").Value = '  +2.36%  '
$ws.Range("D19").Value = '3.200.43'
$ws.Range("E19").Value = '  +1.85%  '
$ws.Range("D20").Value = '''485.36'
$ws.Range("E20").Value = '  +4.75%  '
$ws.Range("D21").Value = '''15.11'
$ws.Range("E21").Value = '  +5.99%  '
$ws.Range("D22").Value = '''0.776'
$ws.Range("E22").Value = '  +6.35%  '
$ws.Range("D23").Value = '''7.94'
$ws.Range("E23").Value = '  +6.65%  '
$ws.Range("E24").Value = '  +7.42%  '
$ws.Range("E25").Value = '  +11.75%  '
$ws.Range("D26").Value = '''83.67'
$ws.Range("E26").Value = '  +2.86%  '
$ws.Range("E27").Value = '  +0.31%  '
$ws.Range("D28").Value = '''9.83'
$ws.Range("E28").Value = '  +7.53%  '
$ws.Range("E29").Value = '  +3.97%  '
$ws.Range("D30").Value = '''2.28'
$ws.Range("E30").Value = '  +3.14%  '
$ws.Range("E31").Value = '  +7.41%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("E33").Value = '  +9.85%  '
$ws.Range("D34").Value = '''28.62'
$ws.Range("E34").Value = '  +6.18%  '
$ws.Range("D35").Value = '0.0₃0898'
$ws.Range("E35").Value = '  +4.95%  '
$ws.Range("D36").Value = '''3.59'
$ws.Range("E36").Value = '  +6.24%  '
$ws.Range("E37").Value = '  +4.79%  '
$ws.Range("E38").Value = '  +5.83%  '
$ws.Range("E39").Value = '  +3.54%  '
$ws.Range("D40").Value = '''474.00'
$ws.Range("E40").Value = '  +7.74%  '
$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").Value = '''9.46'
$ws.Range("E41").Value = '  +7.45%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '''51.68'
$ws.Range("E42").Value = '  +1.85%  '
$ws.Range("E43").Value = '  +9.71%  '
$ws.Range("E44").Value = '  +3.58%  '
$ws.Range("D45").Value = '2.960.65'
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("E46").Value = '  +4.01%  '
$ws.Range("D47").Value = '''38.70'
$ws.Range("E47").Value = '  +5.55%  '
$ws.Range("D48").Value = '''131.69'
$ws.Range("E48").Value = '  +4.74%  '
$ws.Range("E49").Value = '  +8.08%  '
$ws.Range("D50").Value = '''25.63'
$ws.Range("E50").Value = '  +4.95%  '
